# "InputFailureTests" test-case sheet: the "error" column header is being
# renamed to the clearer "expectedError", the column is widened so the new
# header fits, and the sheet's saved selection moves to E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InputFailureTests")
$ws.Activate()

# D1 header: "error" -> "expectedError"
$ws.Range("D1").Value = "expectedError"

# Widen column D to fit the longer header text.
$ws.Columns.Item(4).ColumnWidth = 27.42578125

# Move/save the active cell selection.
$ws.Range("E6").Select()
